$d = $word.ActiveDocument

function Replace-First($old, $new) {
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $rng.Find.Text = $old
    $rng.Find.MatchCase = $true
    $rng.Find.MatchWholeWord = $false
    $rng.Find.Forward = $true
    $rng.Find.Wrap = 0
    $rng.Find.Execute() | Out-Null
    if (-not $rng.Find.Found) {
        throw ("Text not found: " + $old)
    }
    $rng.Text = $new
}

# --- Title ---
Replace-First "Quantum Computing: Unraveling the Enigma" "The Art of Healing: A Journey Through the History of Medicine"

# --- Author name ---
Replace-First "Eleanor Lawson" "Samantha Adams"

# --- Email parts ---
Replace-First "eleanor" "samantha"
Replace-First "lawson@quantum-computing-lab" "adams@school"
Replace-First "org" "edu"

# --- Body paragraph: first 6 sentences (1:1 swap) ---
Replace-First "Step into the realm of quantum computing, where the fabric of reality intertwines with the enigmatic dance of subatomic particles" "From ancient herbal remedies to cutting-edge biotechnology, medicine has a rich and ever-evolving history"
Replace-First "A universe of mind-bending possibilities unfolds before us, challenging our understanding of computation and promising transformative breakthroughs across diverse disciplines" "It is a tale of human ingenuity, resilience, and the pursuit of understanding our bodies and alleviating suffering"
Replace-First "Embark on a journey into this fascinating realm, where the quantum world whispers secrets of untapped potential" "Medicine's journey has been shaped by countless individuals--doctors, scientists, healers, and innovators--who dedicated their lives to pushing the boundaries of knowledge and making a difference in people's lives"

Replace-First "In the heart of a quantum computer, qubits, the quantum counterparts of classical bits, reside in a superposition of states, unlocking the parallel processing of intricate algorithms" "Medicine's origins can be traced back to prehistoric times, with early humans using plants, animal products, and other natural materials to treat illnesses and injuries"
Replace-First "Unlike their classical counterparts, confined to a binary fate, qubits waltz through a symphony of probabilities, traversing multiple paths simultaneously" "As civilizations emerged, so did more formalized systems of medicine, such as those practiced in ancient Egypt, Greece, and China"
Replace-First "This enigmatic ballet of superposition grants quantum computers exponential speed advantages over their classical counterparts, enabling the resolution of previously intractable problems" "These early systems were based on a combination of empirical observations, philosophical beliefs, and religious rituals"

Replace-First "Beyond the realm of theoretical possibilities, quantum computing is poised to revolutionize medicine, materials science, and artificial intelligence" "Over the centuries, medicine underwent profound changes as new ideas and discoveries emerged"
Replace-First "Novel drug discoveries, tailored to individual genetic profiles, hold the promise of personalized healthcare, while quantum algorithms illuminate the path towards previously elusive materials with remarkable properties" "The development of the microscope in the 17th century revolutionized our understanding of the human body and disease"
Replace-First "Artificial intelligence, empowered by quantum computing, embarks on an unprecedented ascent, soaring to new heights of efficiency and accuracy" "The discovery of microorganisms in the 19th century led to the germ theory of disease, which transformed how we approach infection and prevention"

# --- Summary paragraph ---
Replace-First "The enigmatic tapestry of quantum computing is unraveling before our eyes, revealing a world of limitless potential" "Medicine's history is a tale of human endeavor, innovation, and the pursuit of healing"
Replace-First "This mind-bending realm, where superposition dances and qubits pirouette in a quantum waltz, promises transformative breakthroughs across diverse fields, from medicine to materials science to artificial intelligence" "From ancient herbal remedies to modern medical marvels, medicine has undergone profound transformations over time, driven by the dedication of individuals committed to alleviating suffering"
Replace-First "With the dawn of quantum computing, we stand at the precipice of a new era, poised to witness the unfolding of a technological odyssey that will reshape our understanding of reality and redefine the boundaries of human ingenuity" "Despite ongoing challenges, the future of medicine offers hope and promise for a healthier world"
